$d = $word.ActiveDocument

# --- 1. Update the placeholder ID text in the first paragraph, and drop the
#        trailing run that only contained a single space. ---
$d.Content.Find.Execute("**ID__AFFARS_5349_topic_4__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5349_102__ID**", 2)

# --- 2. Give the first paragraph a (5pt-spaced) box border and widen its
#        left indent from 120 to 225 twips (LeftIndent is expressed in
#        points, i.e. twips / 20). ---
$p1 = $d.Paragraphs(1)

$p1.Format.LeftIndent = 225 / 20

$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromRight = 5
